$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1923076923076923
$ws.Range("C2").Value = 0.6153846153846154
$ws.Range("P2").Value = 0.07692307692307693
$ws.Range("S2").Value = 0.1153846153846154
$ws.Range("B3").Value = 0.0625
$ws.Range("C3").Value = 0.0625
$ws.Range("P3").Value = 0.5
$ws.Range("S3").Value = 0.375
$ws.Range("J4").Value = 0.25
$ws.Range("P4").Value = 0.5
$ws.Range("S4").Value = 0.25
$ws.Range("B6").Value = 0.1111111111111111
$ws.Range("F6").Value = 0.05555555555555555
$ws.Range("J6").Value = 0.2222222222222222
$ws.Range("O6").Value = 0.05555555555555555
$ws.Range("Q6").Value = 0.2222222222222222
$ws.Range("S6").Value = 0.3333333333333333
$ws.Range("B7").Value = 0.1612903225806452
$ws.Range("D7").Value = 0.06451612903225806
$ws.Range("F7").Value = 0.03225806451612903
$ws.Range("J7").Value = 0.06451612903225806
$ws.Range("O7").Value = 0.03225806451612903
$ws.Range("Q7").Value = 0.06451612903225806
$ws.Range("R7").Value = 0.03225806451612903
$ws.Range("S7").Value = 0.5483870967741935
$ws.Range("B8").Value = 0.07936507936507936
$ws.Range("F8").Value = 0.07936507936507936
$ws.Range("J8").Value = 0.09523809523809523
$ws.Range("Q8").Value = 0.2063492063492063
$ws.Range("R8").Value = 0.03174603174603174
$ws.Range("S8").Value = 0.5079365079365079
$ws.Range("F9").Value = 0.05263157894736842
$ws.Range("J9").Value = 0.1052631578947368
$ws.Range("O9").Value = 0.05263157894736842
$ws.Range("Q9").Value = 0.2105263157894737
$ws.Range("S9").Value = 0.5789473684210527
$ws.Range("B10").Value = 0.1058823529411765
$ws.Range("D10").Value = 0.02352941176470588
$ws.Range("F10").Value = 0.05882352941176471
$ws.Range("J10").Value = 0.09411764705882353
$ws.Range("Q10").Value = 0.2
$ws.Range("R10").Value = 0.05882352941176471
$ws.Range("S10").Value = 0.4588235294117647
$ws.Range("G11").Value = 0.2727272727272727
$ws.Range("J11").Value = 0.05454545454545454
$ws.Range("K11").Value = 0.3636363636363636
$ws.Range("L11").Value = 0.2727272727272727
$ws.Range("S11").Value = 0.03636363636363636
$ws.Range("G12").Value = 0.7333333333333333
$ws.Range("J12").Value = 0.1333333333333333
$ws.Range("S12").Value = 0.1333333333333333
$ws.Range("G13").Value = 0.7777777777777778
$ws.Range("J13").Value = 0.2222222222222222
$ws.Range("H15").Value = 0.2
$ws.Range("I15").Value = 0.06666666666666667
$ws.Range("J15").Value = 0.3333333333333333
$ws.Range("K15").Value = 0.06666666666666667
$ws.Range("M15").Value = 0.06666666666666667
$ws.Range("O15").Value = 0.06666666666666667
$ws.Range("S15").Value = 0.2
$ws.Range("H16").Value = 0.25
$ws.Range("J16").Value = 0.25
$ws.Range("S16").Value = 0.3333333333333333
$ws.Range("H17").Value = 0.275
$ws.Range("I17").Value = 0.1
$ws.Range("J17").Value = 0.175
$ws.Range("K17").Value = 0.325
$ws.Range("M17").Value = 0.025
$ws.Range("O17").Value = 0.025
$ws.Range("S17").Value = 0.075
$ws.Range("H18").Value = 0.125
$ws.Range("I18").Value = 0.125
$ws.Range("J18").Value = 0.125
$ws.Range("K18").Value = 0.375
$ws.Range("S18").Value = 0.25
$ws.Range("H19").Value = 0.2933333333333333
$ws.Range("I19").Value = 0.09333333333333334
$ws.Range("J19").Value = 0.26
$ws.Range("K19").Value = 0.1133333333333333
$ws.Range("M19").Value = 0.04666666666666667
$ws.Range("O19").Value = 0.05333333333333334
$ws.Range("S19").Value = 0.14
